$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "5.60", "1.00")
# are preserved exactly as typed instead of being coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '40.578.25'
$ws.Range('E2').Value = '  +2.81%  '
$ws.Range('D3').Value = '2.210.16'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '229.12'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').Value = '0.633'
$ws.Range('E6').Value = '  +1.56%  '
$ws.Range('D7').Value = '64.13'
$ws.Range('E7').Value = '  -1.62%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.405'
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('D10').Value = '0.0866'
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '2.538.68'
$ws.Range('E12').Value = '  +1.91%  '
$ws.Range('D13').Value = '15.86'
$ws.Range('E13').Value = '  -1.06%  '
$ws.Range('D14').Value = '22.25'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = '0.825'
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('D16').Value = '5.60'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '2.211.04'
$ws.Range('E17').Value = '  +1.73%  '
$ws.Range('D18').Value = '40.473.44'
$ws.Range('E18').Value = '  +2.52%  '
$ws.Range('D19').Value = '73.92'
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('E20').Value = '  +5.83%  '
$ws.Range('E21').Value = '  -1.02%  '
$ws.Range('D22').Value = '250.07'
$ws.Range('E22').Value = '  +7.66%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('E25').Value = '  -4.92%  '
$ws.Range('D26').Value = '9.69'
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('D27').Value = '173.24'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('D28').Value = '0.141'
$ws.Range('E28').Value = '  +1.46%  '
$ws.Range('D29').Value = '20.36'
$ws.Range('E29').Value = '  +1.33%  '
$ws.Range('E30').Value = '  +2.43%  '
$ws.Range('D31').Value = '2.81'
$ws.Range('E31').Value = '  +1.47%  '
$ws.Range('E32').Value = '  +0.96%  '
$ws.Range('D33').Value = '4.67'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('E34').Value = '  -1.09%  '
$ws.Range('D35').Value = '7.09'
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('E36').Value = '  +1.54%  '
$ws.Range('D37').Value = '3.82'
$ws.Range('E37').Value = '  +5.86%  '
$ws.Range('D38').Value = '2.48'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '4.82'
$ws.Range('E40').Value = '  +11.01%  '
$ws.Range('D41').Value = '0.0232'
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '8.48'
$ws.Range('E42').Value = '  +8.10%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '101.41'
$ws.Range('E43').Value = '  -3.21%  '
$ws.Range('E44').Value = '  +3.59%  '
$ws.Range('D45').Value = '17.52'
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('D46').Value = '1.521.57'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.0933'
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('B48').Value = 'TerraClassic'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D48').Value = '0.000208'
$ws.Range('E48').Value = '  +40.53%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = '1.11'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').Value = '2.81'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').Value = '51.22'
$ws.Range('E51').Value = '  +9.86%  '

# Restore default style on column D (remove the temporary text-number format)
$ws.Range("D2:D51").Style = "Normal"

